$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.15857066666666
$ws.Range("H2").Value = 291.475712
$ws.Range("I2").Value = 0.09537345443416363
$ws.Range("J2").Value = 0.09537345443416365
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 727.4866512122879
$ws.Range("R2").Value = 6547.379860910592
$ws.Range("S2").Value = 0.01320345505455459
$ws.Range("T2").Value = 0.0132034550545546
$ws.Range("G3").Value = 97.15857066666666
$ws.Range("H3").Value = 291.475712
$ws.Range("I3").Value = 0.09537345443416363
$ws.Range("J3").Value = 0.09537345443416365
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 3109.018686630912
$ws.Range("R3").Value = 27981.16817967821
$ws.Range("S3").Value = 0.05642686147477208
$ws.Range("T3").Value = 0.05642686147477208
$ws.Range("G4").Value = 97.15857066666666
$ws.Range("H4").Value = 291.475712
$ws.Range("I4").Value = 0.09537345443416363
$ws.Range("J4").Value = 0.09537345443416365
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 1418.400646550898
$ws.Range("R4").Value = 12765.60581895808
$ws.Range("S4").Value = 0.02574313790483696
$ws.Range("T4").Value = 0.02574313790483697
$ws.Range("G5").Value = 715.8492226666667
$ws.Range("I5").Value = 0.7026967641790764
$ws.Range("J5").Value = 0.7026967641790766
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 5360.008388321832
$ws.Range("R5").Value = 48240.07549489649
$ws.Range("S5").Value = 0.09728100127928167
$ws.Range("T5").Value = 0.09728100127928169
$ws.Range("G6").Value = 715.8492226666667
$ws.Range("I6").Value = 0.7026967641790764
$ws.Range("J6").Value = 0.7026967641790766
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.4157443306037995
$ws.Range("T6").Value = 0.4157443306037995
$ws.Range("G7").Value = 715.8492226666667
$ws.Range("I7").Value = 0.7026967641790764
$ws.Range("J7").Value = 0.7026967641790766
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 10450.55514193262
$ws.Range("R7").Value = 94054.99627739363
$ws.Range("S7").Value = 0.1896714322959953
$ws.Range("T7").Value = 0.1896714322959953
$ws.Range("G8").Value = 205.709325
$ws.Range("H8").Value = 617.127975
$ws.Range("I8").Value = 0.2019297813867599
$ws.Range("J8").Value = 0.2019297813867599
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 1540.27366747515
$ws.Range("R8").Value = 13862.46300727635
$ws.Range("S8").Value = 0.02795506158955979
$ws.Range("T8").Value = 0.02795506158955979
$ws.Range("G9").Value = 205.709325
$ws.Range("H9").Value = 617.127975
$ws.Range("I9").Value = 0.2019297813867599
$ws.Range("J9").Value = 0.2019297813867599
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 6582.5807342661
$ws.Range("R9").Value = 59243.2266083949
$ws.Range("S9").Value = 0.1194699706489836
$ws.Range("T9").Value = 0.1194699706489836
$ws.Range("G10").Value = 205.709325
$ws.Range("H10").Value = 617.127975
$ws.Range("I10").Value = 0.2019297813867599
$ws.Range("J10").Value = 0.2019297813867599
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 3003.113750845375
$ws.Range("R10").Value = 27028.02375760837
$ws.Range("S10").Value = 0.05450474914821644
$ws.Range("T10").Value = 0.05450474914821644
